$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Q4: new year column header (2020), styled like P4 but top-aligned ---
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null
$ws.Range("Q4").Value = 2020
$ws.Range("Q4").VerticalAlignment = -4160

# --- Q5: new data value (21.8), same style as P5 ---
$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5").Value = 21.8

# --- move the active selection to Q9 ---
$ws.Range("Q9").Select() | Out-Null
